$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the scraped cryptos-list refresh: updated Price/Volume(1h) figures,
# plus the coin roster shift (WrappedliquidstakedEther2.0 / Aave dropped,
# Cronos newly appended at the bottom; rows in between shift to match).
$ws.Range("D2").Value = "29.315.66"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "1.841.68"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'238.87"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").Value = "'0.6293"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.07528"
$ws.Range("E8").Value = "  -0.85%  "
$ws.Range("D9").Value = "'0.2943"
$ws.Range("E9").Value = "  -0.78%  "
$ws.Range("D10").Value = "'24.45"
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("D11").Value = "'0.07687"
$ws.Range("D12").Value = "1.855.10"
$ws.Range("E12").Value = "  -6.55%  "
$ws.Range("D13").Value = "'4.969"
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("D14").Value = "'0.6776"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").Value = "'0.00001036"
$ws.Range("E15").Value = "  +4.65%  "
$ws.Range("D16").Value = "'83.06"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Value = "'6.138"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "29.359.05"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'228.45"
$ws.Range("E19").Value = "  -1.23%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'12.41"
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'0.9996"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").Value = "'7.441"
$ws.Range("E22").Value = "  -2.05%  "
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").Value = "'156.26"
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D25").Value = "'0.1392"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'8.350"
$ws.Range("E26").Value = "  -1.30%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'17.61"
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'1.459"
$ws.Range("E28").Value = "  -1.06%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'1.270"
$ws.Range("E29").Value = "  +1.51%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.05643"
$ws.Range("E30").Value = "  -2.81%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.108"
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'4.017"
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("B33").Value = "LidoDAOToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D33").Value = "'1.826"
$ws.Range("E33").Value = "  -2.11%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.154"
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.7087"
$ws.Range("E35").Value = "  -1.16%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.593"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.239.86"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01809"
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.766"
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'6.234"
$ws.Range("E40").Value = "  +2.32%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.9009"
$ws.Range("E41").Value = "  -0.73%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'0.9993"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'101.76"
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'65.44"
$ws.Range("E44").Value = "  -2.90%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "'0.00000000119"
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("D46").Value = "'7.096"
$ws.Range("E46").Value = "  -2.80%  "
$ws.Range("B47").Value = "TheSandbox"
$ws.Range("C47").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D47").Value = "'0.3990"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'8.924"
$ws.Range("E48").Value = "  -2.61%  "
$ws.Range("D49").Value = "'1.671"
$ws.Range("E49").Value = "  -1.47%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.1121"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05714"
$ws.Range("E51").Value = "  -0.67%  "
